$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update blog post title (D9)
$ws.Range("D9").Value = "학위 인증 후기 – 1.어쩌다보니 스위스"

# Row 12: update blog post title (D12) and URL (E12)
$ws.Range("D12").Value = "[혼자 공부하는 머신러닝+딥러닝] 사이킷런 1.1.0 버전 업데이트 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/05/15/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d-%ec%82%ac%ec%9d%b4%ed%82%b7%eb%9f%b0-1-1-0-%eb%b2%84%ec%a0%84-%ec%97%85/"

# Row 46: update blog post title (D46) and URL (E46)
$ws.Range("D46").Value = "[Bioinformatics] 2022년 05월,  미생물 유전체 분석 교육생 모집"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/468"

# Row 51: update blog post title (D51) and URL (E51)
$ws.Range("D51").Value = "[python] seaborn 라이브러리가 제공하는 타이타닉 데이터셋 설명"
$ws.Range("E51").Value = "https://bskyvision.com/1276"
